# Watchlist test cases modified
# Update Runmode (col C) and Results (col D) values on the "Test Cases" sheet,
# and refresh the selection/view state to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Runmode column (C): flip which watchlist cases are enabled -----------
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"
$ws.Range("C8").Value = "N"
$ws.Range("C9").Value = "N"
$ws.Range("C10").Value = "N"
$ws.Range("C11").Value = "N"

# --- Results column (D): updated pass/fail/skip outcomes ------------------
$ws.Range("D2").Value = "FAIL"
$ws.Range("D3").Value = "PASS"
$ws.Range("D4").Value = "PASS"
$ws.Range("D11").Value = "SKIP"

# --- Formatting cleanup: rows 5-11 of column C lose their explicit fill ---
# (the border itself is kept, only the "apply fill" flag goes away)
$ws.Range("C5:C11").Interior.Pattern = -4142
$ws.Range("C6").Borders.Item(8).LineStyle = 1
$ws.Range("C6").Borders.Item(9).LineStyle = 1

# --- Selection / view state -------------------------------------------------
$ws.Activate()
$ws.Range("C2:C4").Select()
